$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")
$ws.Activate() | Out-Null

# --- Row 5: add new formula cell G5 ---
$ws.Range("G5").Formula = "=D5+D7*C6"

# --- Row 13: vertical exchange values updated ---
$ws.Range("B13").Value = 0.066
$ws.Range("C13").Value = 0.066
$ws.Range("D13").Value = 0.002
$ws.Range("E13").Value = "uses 0.5mm sand and 0.06mm fines"
$ws.Range("E13").HorizontalAlignment = -4131

# --- Row 14: alpha coefficient values + number format updated ---
$ws.Range("B14").Value = 0.005
$ws.Range("B14").NumberFormat = "0.00E+00"
$ws.Range("C14").Value = 0.42
$ws.Range("C14").NumberFormat = "0.00"
$ws.Range("D14").Value = 0.16
$ws.Range("D14").NumberFormat = "0.00"

# --- Row 18: hydraulic depth now computed from volume/area ---
$ws.Range("B18").Formula = "=B5/B6"
$ws.Range("C18").Formula = "=C5/C6"
$ws.Range("D18").Formula = "=D5/D6"

# --- Row 19: peak flow velocity updated ---
$ws.Range("B19").Value = 1.5
$ws.Range("C19").Value = 1.5
$ws.Range("D19").Value = 1

# --- Update the active selection to reflect where the author was working ---
$ws.Range("G6").Select() | Out-Null
